# Completed long term outlook view.
# Insert a new intro row on the VARS sheet (sheet1 / first worksheet) at row 4,
# shifting the existing content down, and populate it with the fsr_intro
# variable and its English/French descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VARS")

# Insert a new row before the current row 4 (pushes everything down by one).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the fsr_intro variable and its translations.
$ws.Range("A4").Value = "fsr_intro"
$ws.Range("B4").Value = "The OPBO has prepared fact sheets for the provinces and territories.  Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat."
$ws.Range("C4").Value = "Le BDPB a préparé des fiches d'information pour les provinces et les territoires.  Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat."

# Match the new row's height to the rest of the short, single-line rows.
$ws.Rows.Item(4).RowHeight = 30

# Update the active selection to the newly inserted row, as seen in the file.
$ws.Range("A4").Select()
